$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 62, shifting existing rows 62:165 down to 63:166
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with the new weekly record
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 45100
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 100112001
$ws.Range("G62").Value = "Berenjena"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 180
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = 6444
$ws.Range("N62").Value = "$/caja 60 unidades"
$ws.Range("O62").Value = "Región de Arica y Parinacota"
$ws.Range("P62").Value = 107
$ws.Range("Q62").Value = 60
$ws.Range("R62").Value = "Hortaliza"
